$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.136.54'
$ws.Range('E2').Value = '  -0.48%  '
$ws.Range('D3').Value = '3.307.86'
$ws.Range('E3').Value = '  -1.33%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '188.25'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '556.23'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.33%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.20%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.586'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.71%  '
$ws.Range('D9').Value = '3.298.69'
$ws.Range('E9').Value = '  -1.35%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.185'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.80%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.586'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.93%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '47.42'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.56%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000271'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.11%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.65'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.36%  '
$ws.Range('D15').Value = '3.839.23'
$ws.Range('E15').Value = '  -1.39%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '613.30'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.77%  '
$ws.Range('D17').Value = '66.168.54'
$ws.Range('E17').Value = '  -0.06%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.05'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.04%  '
$ws.Range('E19').Value = '  -0.15%  '
$ws.Range('D20').Value = '3.307.99'
$ws.Range('E20').Value = '  -1.65%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.98'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.48%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.906'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.20%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '18.21'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +8.92%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '102.48'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +6.09%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.96'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.23%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.94'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.90%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.04'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.56%  '
$ws.Range('E28').Value = '  -0.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.63'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.99%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.65'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.29%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '30.24'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.25%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.02'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.82%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.46'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.38%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.06'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.67%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '557.32'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.90%  '
$ws.Range('D36').Value = '3.859.78'
$ws.Range('E36').Value = '  +0.86%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.105'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.06%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '57.34'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.52%  '
$ws.Range('E39').Value = '  +0.09%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.32'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.85%  '
$ws.Range('B41').Value = 'PEPE'
$ws.Range('C41').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D41').Value = '0.0₃0725'
$ws.Range('E41').Value = '  +0.37%  '
$ws.Range('B42').Value = 'InjectiveProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '33.94'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.77%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.73'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.36%  '
$ws.Range('E44').Value = '  +0.72%  '
$ws.Range('E45').Value = '  -3.03%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.21'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -13.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0420'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.95%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.21'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.12%  '
$ws.Range('E49').Value = '  -0.42%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.58'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.99%  '
$ws.Range('E51').Value = '  +0.03%  '
